{"js": "// The commit re-orders the child elements inside <w:rPr> for a set of\n// \"*Tok\" character styles in styles.xml so that <w:b/>/<w:i/> come before\n// <w:color/> (matching the wml.xsd CT_RPr sequence, fixing an\n// OOXMLValidator schema warning). No formatting values actually change -\n// the styles were already bold/italic; only the serialized element order\n// changes. We reproduce that by re-asserting bold/italic (true -> true)\n// on each affected style's font, which causes the run-properties to be\n// rewritten in schema order.\n\nconst styleNames = [\n  \"KeywordTok\",\n  \"ImportTok\",\n  \"CommentTok\",\n  \"DocumentationTok\",\n  \"AnnotationTok\",\n  \"CommentVarTok\",\n  \"ControlFlowTok\",\n  \"InformationTok\",\n  \"WarningTok\",\n  \"AlertTok\",\n  \"ErrorTok\",\n];\n\nconst styles = [];\nfor (const name of styleNames) {\n  const style = context.document.getStyles().getByNameOrNullObject(name);\n  style.load(\"font,nameLocal\");\n  styles.push(style);\n}\nawait context.sync();\n\nfor (const style of styles) {\n  if (style.isNullObject) {\n    continue;\n  }\n  if (style.font.bold) {\n    style.font.bold = true;\n  }\n  if (style.font.italic) {\n    style.font.italic = true;\n  }\n}\nawait context.sync();\n", "ps1": "# The commit re-orders the child elements inside <w:rPr> for a set of\n# \"*Tok\" character styles in styles.xml so that <w:b/>/<w:i/> come before\n# <w:color/> (matching the wml.xsd CT_RPr sequence, fixing an\n# OOXMLValidator schema warning). No formatting values actually change -\n# the styles were already bold/italic; only the serialized element order\n# changes. We reproduce that by re-asserting Bold/Italic (True -> True)\n# on each affected style's Font, which causes the run-properties to be\n# rewritten in schema order.\n\n$d = $word.ActiveDocument\n\n$styleNames = @(\n    \"KeywordTok\",\n    \"ImportTok\",\n    \"CommentTok\",\n    \"DocumentationTok\",\n    \"AnnotationTok\",\n    \"CommentVarTok\",\n    \"ControlFlowTok\",\n    \"InformationTok\",\n    \"WarningTok\",\n    \"AlertTok\",\n    \"ErrorTok\"\n)\n\nforeach ($name in $styleNames) {\n    $s = $d.Styles($name)\n    if ($s.Font.Bold) {\n        $s.Font.Bold = 1\n    }\n    if ($s.Font.Italic) {\n        $s.Font.Italic = 1\n    }\n}\n"}
